# "Generate Report for Handback"
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps for the first
# handback entry (02c9d5e4-...) on both the "zh-cn" and "de-de" report
# sheets. Rows 2 and 4 on each sheet share the same underlying values,
# so both are updated to keep the workbook internally consistent.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-18 07:00:14"
$wsZh.Range("E4").Value = "2016-03-18 07:00:14"
$wsZh.Range("H2").Value = "2016-03-18 07:00:32"
$wsZh.Range("H4").Value = "2016-03-18 07:00:32"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-18 07:00:17"
$wsDe.Range("E4").Value = "2016-03-18 07:00:17"
$wsDe.Range("H2").Value = "2016-03-18 07:00:42"
$wsDe.Range("H4").Value = "2016-03-18 07:00:42"
